$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'CreatedAt: 2025-12-28T18:06:43'
$ws.Range("V4").Value = 91.79000000000001
$ws.Range("W4").Value = 69.31999999999999
$ws.Range("X4").Value = 43.52
$ws.Range("Y4").Value = 44.22
$ws.Range("Z4").Value = 43.93
$ws.Range("V6").Value = -1.74
$ws.Range("W6").Value = -0.9
$ws.Range("X6").Value = -0.65
$ws.Range("Y6").Value = -0.09
$ws.Range("Z6").Value = 0.22
$ws.Range("V9").Value = 92.42
$ws.Range("W9").Value = 69.66
$ws.Range("X9").Value = 43.91
$ws.Range("Y9").Value = 45.07
$ws.Range("Z9").Value = 45.2
$ws.Range("V11").Value = -1.11
$ws.Range("W11").Value = -0.5600000000000001
$ws.Range("X11").Value = -0.26
$ws.Range("Y11").Value = 0.77
$ws.Range("Z11").Value = 1.49
$ws.Range("V14").Value = 92.42
$ws.Range("W14").Value = 69.66
$ws.Range("X14").Value = 43.91
$ws.Range("Y14").Value = 45.07
$ws.Range("Z14").Value = 45.2
$ws.Range("V16").Value = -1.11
$ws.Range("W16").Value = -0.5600000000000001
$ws.Range("X16").Value = -0.26
$ws.Range("Y16").Value = 0.77
$ws.Range("Z16").Value = 1.49
$ws.Range("V19").Value = 45.91
$ws.Range("W19").Value = 45.12
$ws.Range("X19").Value = 42.97
$ws.Range("Y19").Value = 43.69
$ws.Range("Z19").Value = 43.49
$ws.Range("V20").Value = -44.55
$ws.Range("W20").Value = -23.18
$ws.Range("V21").Value = -3.08
$ws.Range("W21").Value = -1.91
$ws.Range("X21").Value = -1.2
$ws.Range("Y21").Value = -0.61
$ws.Range("Z21").Value = -0.22
$ws.Range("V24").Value = 90.45999999999999
$ws.Range("W24").Value = 68.31
$ws.Range("X24").Value = 42.97
$ws.Range("Y24").Value = 43.69
$ws.Range("Z24").Value = 43.49
$ws.Range("V26").Value = -3.08
$ws.Range("W26").Value = -1.91
$ws.Range("X26").Value = -1.2
$ws.Range("Y26").Value = -0.61
$ws.Range("Z26").Value = -0.22
$ws.Range("V29").Value = 44.19
$ws.Range("W29").Value = 43.82
$ws.Range("X29").Value = 42.23
$ws.Range("Y29").Value = 42.89
$ws.Range("Z29").Value = 42.85
$ws.Range("V30").Value = -44.55
$ws.Range("W30").Value = -23.18
$ws.Range("V31").Value = -4.79
$ws.Range("X31").Value = -1.94
$ws.Range("Y31").Value = -1.42
$ws.Range("Z31").Value = -0.86
$ws.Range("V34").Value = 95.25
$ws.Range("W34").Value = 71.95
$ws.Range("X34").Value = 45.08
$ws.Range("Y34").Value = 46.79
$ws.Range("Z34").Value = 46.95
$ws.Range("V36").Value = 1.71
$ws.Range("W36").Value = 1.73
$ws.Range("X36").Value = 0.9
$ws.Range("Y36").Value = 2.48
$ws.Range("Z36").Value = 3.24
$ws.Range("Z37").Value = 0
$ws.Range("V39").Value = 91.79000000000001
$ws.Range("W39").Value = 69.31999999999999
$ws.Range("X39").Value = 43.52
$ws.Range("Y39").Value = 44.22
$ws.Range("Z39").Value = 43.93
$ws.Range("V41").Value = -1.74
$ws.Range("W41").Value = -0.9
$ws.Range("X41").Value = -0.65
$ws.Range("Y41").Value = -0.09
$ws.Range("Z41").Value = 0.22
$ws.Range("V44").Value = 101.12
$ws.Range("W44").Value = 75.83
$ws.Range("X44").Value = 47.55
$ws.Range("Y44").Value = 47.59
$ws.Range("Z44").Value = 47
$ws.Range("V46").Value = 7.58
$ws.Range("W46").Value = 5.61
$ws.Range("X46").Value = 3.38
$ws.Range("Y46").Value = 3.28
$ws.Range("Z46").Value = 3.29
$ws.Range("V49").Value = 98.98
$ws.Range("X49").Value = 47.14
$ws.Range("Y49").Value = 47.74
$ws.Range("Z49").Value = 46.55
$ws.Range("V51").Value = 5.44
$ws.Range("X51").Value = 2.97
$ws.Range("Y51").Value = 3.44
$ws.Range("Z51").Value = 2.84
$ws.Range("V54").Value = 94.38
$ws.Range("X54").Value = 45.82
$ws.Range("Y54").Value = 46.44
$ws.Range("Z54").Value = 45.57
$ws.Range("V56").Value = 0.85
$ws.Range("X56").Value = 1.65
$ws.Range("Y56").Value = 2.14
$ws.Range("Z56").Value = 1.87
$ws.Range("V59").Value = 99.95
$ws.Range("W59").Value = 75.06999999999999
$ws.Range("X59").Value = 47.07
$ws.Range("Y59").Value = 47.16
$ws.Range("Z59").Value = 46.5
$ws.Range("V61").Value = 6.41
$ws.Range("W61").Value = 4.85
$ws.Range("X61").Value = 2.9
$ws.Range("Y61").Value = 2.86
$ws.Range("Z61").Value = 2.79
$ws.Range("V64").Value = 99.95
$ws.Range("W64").Value = 75.06999999999999
$ws.Range("X64").Value = 47.07
$ws.Range("Y64").Value = 47.16
$ws.Range("Z64").Value = 46.5
$ws.Range("V66").Value = 6.41
$ws.Range("W66").Value = 4.85
$ws.Range("X66").Value = 2.9
$ws.Range("Y66").Value = 2.86
$ws.Range("Z66").Value = 2.79
$ws.Range("W69").Value = 270
$ws.Range("X69").Value = 47.98
$ws.Range("Y69").Value = 48.28
$ws.Range("Z69").Value = 47.3
$ws.Range("V71").Value = 8.130000000000001
$ws.Range("X71").Value = 3.69
$ws.Range("Y71").Value = 3.7
$ws.Range("Z71").Value = 3.59
$ws.Range("V72").Value = 73.33
$ws.Range("W72").Value = 193.68
$ws.Range("X72").Value = 0.12
$ws.Range("Y72").Value = 0.28
$ws.Range("V74").Value = 99.08
$ws.Range("X74").Value = 46.79
$ws.Range("Y74").Value = 46.83
$ws.Range("Z74").Value = 46.15
$ws.Range("V76").Value = 5.55
$ws.Range("X76").Value = 2.62
$ws.Range("Y76").Value = 2.53
$ws.Range("Z76").Value = 2.45
$ws.Range("V79").Value = 99.68000000000001
$ws.Range("W79").Value = 75.12
$ws.Range("X79").Value = 47.12
$ws.Range("Y79").Value = 47.18
$ws.Range("Z79").Value = 46.51
$ws.Range("V81").Value = 6.14
$ws.Range("W81").Value = 4.9
$ws.Range("X81").Value = 2.94
$ws.Range("Y81").Value = 2.88
$ws.Range("Z81").Value = 2.8
$ws.Range("V84").Value = 93.06999999999999
$ws.Range("X84").Value = 46.01
$ws.Range("Y84").Value = 46.59
$ws.Range("Z84").Value = 45.48
$ws.Range("V86").Value = -0.47
$ws.Range("X86").Value = 1.84
$ws.Range("Y86").Value = 2.28
$ws.Range("Z86").Value = 1.77
$ws.Range("V89").Value = 88.73999999999999
$ws.Range("W89").Value = 66.94
$ws.Range("X89").Value = 42.23
$ws.Range("Y89").Value = 42.89
$ws.Range("Z89").Value = 42.85
$ws.Range("V91").Value = -4.79
$ws.Range("W91").Value = -3.28
$ws.Range("X91").Value = -1.94
$ws.Range("Y91").Value = -1.42
$ws.Range("Z91").Value = -0.86
